$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "SW"
$ws.Range("B5").Value = "SW=HC=LS=R=SS"
$ws.Range("A8").Value = "true=Wisdom=classes/cleric/SpellSlotsCleric.xlsx=classes/cleric/SpellListCleric.xlsx=P"
